$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '50.587.47'
$ws.Range("E2").Value = '  -1.35%  '

# Row 3
$ws.Range("D3").Value = '2.908.68'
$ws.Range("E3").Value = '  -2.36%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '371.97'
$ws.Range("E5").Value = '  -3.10%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.71'
$ws.Range("E6").Value = '  -3.20%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.532'
$ws.Range("E7").Value = '  -1.41%  '

# Row 8
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.572'
$ws.Range("E9").Value = '  -3.29%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.26'
$ws.Range("E10").Value = '  -3.76%  '

# Row 11
$ws.Range("E11").Value = '  -0.93%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0840'
$ws.Range("E12").Value = '  -0.11%  '

# Row 13
$ws.Range("D13").Value = '3.358.88'
$ws.Range("E13").Value = '  -2.65%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.78'
$ws.Range("E14").Value = '  -2.00%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.44'
$ws.Range("E15").Value = '  -0.43%  '

# Row 16
$ws.Range("B16").Value = 'Uniswap'
$ws.Range("C16").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '11.84'
$ws.Range("E16").Value = '  +60.58%  '

# Row 17
$ws.Range("D17").Value = '2.906.86'
$ws.Range("E17").Value = '  -2.53%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.976'
$ws.Range("E18").Value = '  -0.92%  '

# Row 19
$ws.Range("D19").Value = '50.575.54'
$ws.Range("E19").Value = '  -1.33%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.97'
$ws.Range("E20").Value = '  -8.58%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.05'
$ws.Range("E21").Value = '  -5.54%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0936'
$ws.Range("E22").Value = '  -2.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.83'
$ws.Range("E23").Value = '  +0.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '263.50'
$ws.Range("E24").Value = '  +0.72%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.08'
$ws.Range("E25").Value = '  +6.58%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.71'
$ws.Range("E26").Value = '  -5.42%  '

# Row 27
$ws.Range("E27").Value = '  -0.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.04'
$ws.Range("E28").Value = '  -5.66%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.13'
$ws.Range("E29").Value = '  -2.40%  '

# Row 30
$ws.Range("E30").Value = '  -4.44%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.107'
$ws.Range("E31").Value = '  -6.05%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.80'
$ws.Range("E32").Value = '  -0.34%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.24'
$ws.Range("E33").Value = '  -1.38%  '

# Row 34
$ws.Range("E34").Value = '  -0.92%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '32.62'
$ws.Range("E35").Value = '  -4.91%  '

# Row 36
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  -0.08%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0424'
$ws.Range("E37").Value = '  -4.95%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.03'
$ws.Range("E38").Value = '  +1.46%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.114'
$ws.Range("E39").Value = '  -1.21%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.10'
$ws.Range("E40").Value = '  -5.21%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.75'
$ws.Range("E41").Value = '  -2.05%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.38'
$ws.Range("E42").Value = '  -7.86%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '118.68'
$ws.Range("E43").Value = '  -2.79%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.56'
$ws.Range("E44").Value = '  -4.73%  '

# Row 45
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.03'
$ws.Range("E45").Value = '  -2.31%  '

# Row 46
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.31'
$ws.Range("E46").Value = '  +2.43%  '

# Row 47
$ws.Range("E47").Value = '  -0.82%  '

# Row 48
$ws.Range("D48").Value = '1.977.57'
$ws.Range("E48").Value = '  -2.24%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.254'
$ws.Range("E49").Value = '  -7.64%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0312'
$ws.Range("E50").Value = '  -5.05%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.20'
$ws.Range("E51").Value = '  +2.02%  '
